$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "297.42"
$ws.Range("E2").NumberFormat = "@"
$ws.Range("E2").Value = "2.58%"

$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "41.54"
$ws.Range("E3").NumberFormat = "@"
$ws.Range("E3").Value = "3.54%"

$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "5.028"
$ws.Range("E4").NumberFormat = "@"
$ws.Range("E4").Value = "-0.36%"

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "0.07552"
$ws.Range("E5").NumberFormat = "@"
$ws.Range("E5").Value = "3.49%"

$ws.Range("B6").Value = "FTXToken"
$ws.Range("C6").Value = "https://coinranking.com/coin/NfeOYfNcl+ftxtoken-ftt"
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "1.601"
$ws.Range("E6").NumberFormat = "@"
$ws.Range("E6").Value = "3.14%"

$ws.Range("B7").Value = "MXToken"
$ws.Range("C7").Value = "https://coinranking.com/coin/QUC5kVAxSoB-+mxtoken-mx"
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.9280"
$ws.Range("E7").NumberFormat = "@"
$ws.Range("E7").Value = "1.04%"

$ws.Range("B8").Value = "BTSEToken"
$ws.Range("C8").Value = "https://coinranking.com/coin/EOSL_JJKNMEr+btsetoken-btse"
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "2.404"
$ws.Range("E8").NumberFormat = "@"
$ws.Range("E8").Value = "0.30%"

$ws.Range("B9").Value = "LiechtensteinCryptoassetsExchange"
$ws.Range("C9").Value = "https://coinranking.com/coin/v4IW9oaF+liechtensteincryptoassetsexchange-lcx"
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.1199"
$ws.Range("E9").NumberFormat = "@"
$ws.Range("E9").Value = "3.11%"

$ws.Range("B10").Value = "WazirX"
$ws.Range("C10").Value = "https://coinranking.com/coin/6QK-8hUZ+wazirx-wrx"
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "0.1841"
$ws.Range("E10").NumberFormat = "@"
$ws.Range("E10").Value = "6.62%"

$ws.Range("B11").Value = "MandalaExchangeToken"
$ws.Range("C11").Value = "https://coinranking.com/coin/lviNIbma2Xuqs+mandalaexchangetoken-mdx"
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.08923"
$ws.Range("E11").NumberFormat = "@"
$ws.Range("E11").Value = "2.92%"

$ws.Range("B12").Value = "BitrueCoin"
$ws.Range("C12").Value = "https://coinranking.com/coin/SLYjzF4ty+bitruecoin-btr"
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "0.04037"
$ws.Range("E12").NumberFormat = "@"
$ws.Range("E12").Value = "-3.13%"

$ws.Range("B13").Value = "BitMartToken"
$ws.Range("C13").Value = "https://coinranking.com/coin/6uzcPMFgWUJNH+bitmarttoken-bmx"
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "0.1054"
$ws.Range("E13").NumberFormat = "@"
$ws.Range("E13").Value = "0.09%"

$ws.Range("B14").Value = "BitForexToken"
$ws.Range("C14").Value = "https://coinranking.com/coin/2nh5ugplNocUp+bitforextoken-bf"
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "0.001281"
$ws.Range("E14").NumberFormat = "@"
$ws.Range("E14").Value = "0.88%"

$ws.Range("B15").Value = "TigerCash"
$ws.Range("C15").Value = "https://coinranking.com/coin/6hIn06L2+tigercash-tch"
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "0.005816"
$ws.Range("E15").NumberFormat = "@"
$ws.Range("E15").Value = "0.66%"

$ws.Range("B16").Value = "LEO"
$ws.Range("C16").Value = "https://coinranking.com/coin/mqtUpyBxu8O8+leo-leo"
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "3.339"
$ws.Range("E16").NumberFormat = "@"
$ws.Range("E16").Value = "-1.55%"

$ws.Range("B17").Value = "GateToken"
$ws.Range("C17").Value = "https://coinranking.com/coin/t7m8DZVyMsAu+gatetoken-gt"
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "4.375"
$ws.Range("E17").NumberFormat = "@"
$ws.Range("E17").Value = "2.28%"

$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "7.969"
$ws.Range("E19").NumberFormat = "@"
$ws.Range("E19").Value = "1.60%"

$ws.Range("E20").NumberFormat = "@"
$ws.Range("E20").Value = "5.01%"

$ws.Range("E21").NumberFormat = "@"
$ws.Range("E21").Value = "4.00%"

$ws.Range("E22").NumberFormat = "@"
$ws.Range("E22").Value = "4.88%"

$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "0.001265"
$ws.Range("E23").NumberFormat = "@"
$ws.Range("E23").Value = "-0.33%"

$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "0.004152"
$ws.Range("E24").NumberFormat = "@"
$ws.Range("E24").Value = "6.96%"

$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "0.0001230"
$ws.Range("E25").NumberFormat = "@"
$ws.Range("E25").Value = "-4.00%"

$ws.Range("E26").NumberFormat = "@"
$ws.Range("E26").Value = "-0.02%"

$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "0.02417"
$ws.Range("E38").NumberFormat = "@"
$ws.Range("E38").Value = "4.09%"

$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "0.05209"
$ws.Range("E39").NumberFormat = "@"
$ws.Range("E39").Value = "4.91%"

$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "0.006508"
$ws.Range("E40").NumberFormat = "@"
$ws.Range("E40").Value = "-4.78%"

$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "0.007800"
$ws.Range("E41").NumberFormat = "@"
$ws.Range("E41").Value = "1.81%"

$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "0.1328"

$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "0.007552"
$ws.Range("E43").NumberFormat = "@"
$ws.Range("E43").Value = "2.71%"

$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "0.007833"
$ws.Range("E44").NumberFormat = "@"
$ws.Range("E44").Value = "10.87%"

$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "0.3222"
$ws.Range("E45").NumberFormat = "@"
$ws.Range("E45").Value = "11.55%"

$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "0.00006783"

$ws.Range("E47").NumberFormat = "@"
$ws.Range("E47").Value = "-0.10%"

$ws.Range("E48").NumberFormat = "@"
$ws.Range("E48").Value = "-0.08%"

$ws.Range("E49").NumberFormat = "@"
$ws.Range("E49").Value = "441.87%"

$ws.Range("E50").NumberFormat = "@"
$ws.Range("E50").Value = "-0.10%"

$ws.Range("E51").NumberFormat = "@"
$ws.Range("E51").Value = "-0.10%"
